# CDS Input file updates
# Replace the ParticipantsTab Cypher query (cell B2 on the "startup" sheet)
# with the new/updated query text, then refresh the row height and the
# saved cursor/selection position to match the resulting layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$newParticipantsQuery = @'
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE f.file_type in ['IDAT']
with p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN
coalesce(p.participant_id,'') as `Participant ID`,
coalesce(s.study_name, '') as `Study Name`,
coalesce(s.phs_accession,'') as `Accession`,
coalesce(p.gender,'') as `Gender`,
coalesce(apoc.text.join(samp, ','), '') as `Samples`
ORDER BY p.participant_id LIMIT 100
'@

# Cell B2 holds the query for the "ParticipantsTab" row.
$ws.Range("B2").Value = $newParticipantsQuery

# The replacement text wraps onto more lines, so the row needs to grow to
# keep the whole query visible.
$ws.Rows.Item(2).RowHeight = 279

# Update the saved cursor position / scrolled view.
$ws.Range("B5").Select()
